$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -207
$ws.Range("N10").Value = -1086
$ws.Range("H113").Value = 1807.8846
$ws.Range("I113").Value = 1587.174
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1587.174
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 1666.826
$ws.Range("N113").Value = -10008
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -87
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9407
$ws.Range("H58").Value = 1978.24
$ws.Range("I58").Value = 1908
$ws.Range("K58").Value = 1908
$ws.Range("M58").Value = -1705
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 27673.611
$ws.Range("J68").Value = 30007.812
$ws.Range("L68").Value = 30007.812
$ws.Range("N68").Value = -31505.812
$ws.Range("H71").Value = 27673.611
$ws.Range("J71").Value = 30007.812
$ws.Range("L71").Value = 90023.436
$ws.Range("N71").Value = -97511.436
$ws.Range("H99").Value = 1758
$ws.Range("I99").Value = 1366.4445
$ws.Range("J99").Value = 2462.8
$ws.Range("K99").Value = 1366.4445
$ws.Range("L99").Value = 2462.8
$ws.Range("M99").Value = 131.5554999999999
$ws.Range("N99").Value = -5458.8
$ws.Range("H117").Value = 73200
$ws.Range("J117").Value = 73200
$ws.Range("L117").Value = 73200
$ws.Range("N117").Value = -82378
$ws.Range("H126").Value = 1758
$ws.Range("I126").Value = 1366.4445
$ws.Range("J126").Value = 2462.8
$ws.Range("K126").Value = 4099.333500000001
$ws.Range("L126").Value = 7388.400000000001
$ws.Range("M126").Value = -1629.333500000001
$ws.Range("N126").Value = -12328.4
$ws.Range("H136").Value = 1978.24
$ws.Range("I136").Value = 1908
$ws.Range("K136").Value = 5724
$ws.Range("M136").Value = -3174

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1942.8235
$ws.Range("I3").Value = 1442
$ws.Range("J3").Value = 2506.25
$ws.Range("K3").Value = 4326
$ws.Range("L3").Value = 7518.75
$ws.Range("M3").Value = -4214
$ws.Range("N3").Value = -7742.75
$ws.Range("H7").Value = 88.666664
$ws.Range("I7").Value = 83
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 249
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -137
$ws.Range("N7").Value = -524
$ws.Range("H17").Value = 720.63635
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 1224.5
$ws.Range("K17").Value = 348
$ws.Range("L17").Value = 3673.5
$ws.Range("M17").Value = -179
$ws.Range("N17").Value = -4011.5
$ws.Range("H34").Value = 1375.9286
$ws.Range("J34").Value = 1743.5
$ws.Range("L34").Value = 5230.5
$ws.Range("N34").Value = -5398.5
$ws.Range("H39").Value = 1308.5385
$ws.Range("I39").Value = 900
$ws.Range("J39").Value = 1342.5834
$ws.Range("K39").Value = 2700
$ws.Range("L39").Value = 4027.7502
$ws.Range("M39").Value = -2406
$ws.Range("N39").Value = -4615.7502
$ws.Range("H55").Value = 2855
$ws.Range("I55").Value = 1854
$ws.Range("J55").Value = 2998
$ws.Range("K55").Value = 5562
$ws.Range("L55").Value = 8994
$ws.Range("M55").Value = -5385
$ws.Range("N55").Value = -9348
$ws.Range("H64").Value = 2415.4285
$ws.Range("J64").Value = 2482.6667
$ws.Range("L64").Value = 7448.000100000001
$ws.Range("N64").Value = -7988.000100000001
$ws.Range("H67").Value = 2415.4285
$ws.Range("J67").Value = 2482.6667
$ws.Range("L67").Value = 7448.000100000001
$ws.Range("N67").Value = -9320.000100000001
$ws.Range("H134").Value = 3261.8235
$ws.Range("I134").Value = 1328.9524
$ws.Range("J134").Value = 4614.8335
$ws.Range("K134").Value = 3986.857199999999
$ws.Range("L134").Value = 13844.5005
$ws.Range("M134").Value = 1083.142800000001
$ws.Range("N134").Value = -23984.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H43").Value = 2500
$ws.Range("I43").Value = 2500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -2349
$ws.Range("N43").ClearContents()
$ws.Range("H92").Value = 3125.5
$ws.Range("J92").Value = 3125.5
$ws.Range("L92").Value = 3125.5
$ws.Range("N92").Value = -6869.5
$ws.Range("H102").Value = 1232.25
$ws.Range("I102").Value = 1093
$ws.Range("J102").Value = 1650
$ws.Range("K102").Value = 1093
$ws.Range("L102").Value = 1650
$ws.Range("M102").Value = 529
$ws.Range("N102").Value = -4894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2705.5293
$ws.Range("I40").Value = 2383.7778
$ws.Range("J40").Value = 3067.5
$ws.Range("K40").Value = 2383.7778
$ws.Range("L40").Value = 3067.5
$ws.Range("M40").Value = -2247.7778
$ws.Range("N40").Value = -3339.5
$ws.Range("H55").Value = 486.57144
$ws.Range("I55").Value = 513
$ws.Range("J55").Value = 451.33334
$ws.Range("K55").Value = 513
$ws.Range("L55").Value = 451.33334
$ws.Range("M55").Value = -340
$ws.Range("N55").Value = -797.33334
$ws.Range("H122").Value = 3605.7896
$ws.Range("I122").Value = 3278.077
$ws.Range("J122").Value = 4315.8335
$ws.Range("K122").Value = 9834.231
$ws.Range("L122").Value = 12947.5005
$ws.Range("M122").Value = -7384.231
$ws.Range("N122").Value = -17847.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30346
$ws.Range("H122").Value = 286984.5
$ws.Range("I122").Value = 527228.5
$ws.Range("J122").Value = 1694.6875
$ws.Range("K122").Value = 1581685.5
$ws.Range("L122").Value = 5084.0625
$ws.Range("M122").Value = -1579235.5
$ws.Range("N122").Value = -9984.0625

